# Generate Report for Handback
#
# The localization-status workbook tracks, per target locale sheet
# (zh-cn / de-de), the handoff of source files for translation. This
# script records that the translations have now been handed BACK:
#   - Column C ("Status") flips from "Ready for handoff" to
#     "Handed back: in sync with en-US".
#   - Column F ("Latest Target File") and column G ("Latest Handback
#     File") get populated (they were blank while awaiting handback),
#     mirroring the md/xlf file names already linked in columns A/D,
#     each as a real hyperlink (consistent with how A/B/D are already
#     rendered as hyperlinks).
#   - Column H ("Latest Handback DateTime") gets stamped with the
#     actual handback time instead of the "0001-01-01 00:00:00"
#     placeholder.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

function Update-LocaleSheet($SheetName, $MdBaseUrl, $XlfBaseUrl, $Row2Md, $Row2Xlf, $Row3Md, $Row3Xlf, $Row2Handback, $Row3Handback) {

    $ws = $wb.Worksheets.Item($SheetName)

    # -- Row 2 (first file) --
    $ws.Range("C2").Value = $statusText
    $ws.Hyperlinks.Add($ws.Range("F2"), ($MdBaseUrl + $Row2Md), "", "", $Row2Md)
    $ws.Hyperlinks.Add($ws.Range("G2"), ($XlfBaseUrl + $Row2Xlf), "", "", $Row2Xlf)
    $ws.Range("H2").Value = $Row2Handback

    # -- Row 3 (second file) --
    $ws.Range("C3").Value = $statusText
    $ws.Hyperlinks.Add($ws.Range("F3"), ($MdBaseUrl + $Row3Md), "", "", $Row3Md)
    $ws.Hyperlinks.Add($ws.Range("G3"), ($XlfBaseUrl + $Row3Xlf), "", "", $Row3Xlf)
    $ws.Range("H3").Value = $Row3Handback
}

$mdBaseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/7e8699305d652c36cab50265062612bb2ad2b01d/e2e/"

# zh-cn sheet
Update-LocaleSheet `
    "zh-cn" `
    $mdBaseUrl `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f4e94bdf32f641b0affeee2806a8178a125ed84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/" `
    "387d23b5-4926-4679-8d3b-d6ac8679c093.md" `
    "387d23b5-4926-4679-8d3b-d6ac8679c093.4ed0716f53e5e091ee08371bdb9850b45265619b.zh-cn.xlf" `
    "6c84ca05-1e47-44dd-99ba-8800c17c163e.md" `
    "6c84ca05-1e47-44dd-99ba-8800c17c163e.e301d1b696427c9851845275801f72f230385adf.zh-cn.xlf" `
    "2016-03-19 14:16:18" `
    "2016-03-19 14:16:18"

# de-de sheet
Update-LocaleSheet `
    "de-de" `
    $mdBaseUrl `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1ae988720ea9c967315656f129966479910c4052/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/" `
    "387d23b5-4926-4679-8d3b-d6ac8679c093.md" `
    "387d23b5-4926-4679-8d3b-d6ac8679c093.4ed0716f53e5e091ee08371bdb9850b45265619b.de-de.xlf" `
    "6c84ca05-1e47-44dd-99ba-8800c17c163e.md" `
    "6c84ca05-1e47-44dd-99ba-8800c17c163e.e301d1b696427c9851845275801f72f230385adf.de-de.xlf" `
    "2016-03-19 14:16:22" `
    "2016-03-19 14:16:22"

Write-Host "Handback report generated."
